# Pre-Champs Updates (LMR to FMB)
$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Update the cached "datetimeFigureOut" footer field from 3/28/2017 to
#    4/17/2017 everywhere it appears: once on the slide master, and once on
#    each of the slide layouts' own "Date Placeholder" shape.
# ---------------------------------------------------------------------------
$oldDate = "3/28/2017"
$newDate = "4/17/2017"

$master = $p.SlideMaster

for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shp = $master.Shapes.Item($i)
    if ($shp.HasTextFrame -eq -1) {
        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -eq $oldDate) {
            $tr.Text = $newDate
        }
    }
}

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $shp = $layout.Shapes.Item($i)
        if ($shp.HasTextFrame -eq -1) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# ---------------------------------------------------------------------------
# 2) Rename the "Gears Scored LMR" (Left/Mid/Right) labels on slide 2 to
#    "Gears Scored FMB" (Feeder/Mid/Boiler), updating the underlying field
#    names from left/right to fdr/boi for each of the six robot panels.
# ---------------------------------------------------------------------------
function Update-GearsLabel($shape) {
    if ($shape.HasTextFrame -eq -1) {
        $tr = $shape.TextFrame.TextRange
        $lineCount = $tr.Lines().Count
        for ($j = 1; $j -le $lineCount; $j++) {
            $ln = $tr.Lines($j)
            if ($ln.Text -like "Gears Scored LMR:*") {
                $newText = $ln.Text.Replace("Gears Scored LMR:", "Gears Scored FMB:")
                $newText = $newText.Replace("_left_", "_fdr_")
                $newText = $newText.Replace("_right_", "_boi_")
                $ln.Text = $newText
            }
        }
    }
    if ($shape.Type -eq 6) {
        for ($k = 1; $k -le $shape.GroupItems.Count; $k++) {
            Update-GearsLabel $shape.GroupItems.Item($k)
        }
    }
}

$slide2 = $p.Slides.Item(2)
for ($i = 1; $i -le $slide2.Shapes.Count; $i++) {
    Update-GearsLabel $slide2.Shapes.Item($i)
}
